$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = [double]"0.06349258695500602"
$ws.Range("E2").Value = [double]"0.03777348563721421"
$ws.Range("F2").Value = [double]"0.0684526912477503"
$ws.Range("G2").Value = [double]"-0.03158809536073708"
$ws.Range("H2").Value = [double]"0.003923958723663769"
$ws.Range("I2").Value = [double]"0.03903941388948123"
$ws.Range("J2").Value = [double]"0.1721239418668228"
$ws.Range("K2").Value = [double]"0.198667186370348"
$ws.Range("L2").Value = [double]"-0.05612525113176752"
$ws.Range("M2").Value = [double]"0.1130601880950822"
$ws.Range("D3").Value = [double]"0.04608398681596151"
$ws.Range("E3").Value = [double]"0.02189183347268834"
$ws.Range("F3").Value = [double]"0.02012647564497602"
$ws.Range("G3").Value = [double]"0.1693011165915949"
$ws.Range("H3").Value = [double]"0.08724973029280536"
$ws.Range("I3").Value = [double]"-0.05039689115876729"
$ws.Range("J3").Value = [double]"-0.08104252506935648"
$ws.Range("K3").Value = [double]"-0.05592584816530667"
$ws.Range("L3").Value = [double]"0.09867565890553331"
$ws.Range("M3").Value = [double]"0.1012401029327586"
$ws.Range("B4").Value = [double]"0.06349258695500602"
$ws.Range("C4").Value = [double]"0.04608398681596151"
$ws.Range("E4").Value = [double]"0.5261472653192818"
$ws.Range("F4").Value = [double]"0.5543526031480206"
$ws.Range("G4").Value = [double]"0.5482195976645435"
$ws.Range("H4").Value = [double]"0.3021693211650757"
$ws.Range("I4").Value = [double]"0.8662642689568739"
$ws.Range("J4").Value = [double]"0.5044893085154387"
$ws.Range("K4").Value = [double]"0.4317035091799344"
$ws.Range("L4").Value = [double]"0.5059907001350852"
$ws.Range("M4").Value = [double]"0.3384409619207953"
$ws.Range("N4").Value = [double]"7.426571720754092e-16"
$ws.Range("O4").Value = [double]"0.02970424846317648"
$ws.Range("Q4").Value = [double]"4.835907167002664e-16"
$ws.Range("R4").Value = [double]"-0.1987608672475611"
$ws.Range("S4").Value = [double]"0.4302309024326345"
$ws.Range("T4").Value = [double]"0.2205306032889259"
$ws.Range("U4").Value = [double]"0.2244033592937769"
$ws.Range("B5").Value = [double]"0.03777348563721421"
$ws.Range("C5").Value = [double]"0.02189183347268834"
$ws.Range("D5").Value = [double]"0.5261472653192818"
$ws.Range("F5").Value = [double]"0.9818426010967805"
$ws.Range("G5").Value = [double]"0.1320648805310834"
$ws.Range("H5").Value = [double]"0.6773768707777231"
$ws.Range("I5").Value = [double]"0.6355324873451915"
$ws.Range("J5").Value = [double]"0.7614718058408736"
$ws.Range("K5").Value = [double]"0.6007500003764545"
$ws.Range("L5").Value = [double]"0.2210786576314376"
$ws.Range("M5").Value = [double]"0.4575882349190203"
$ws.Range("N5").Value = [double]"-1.079523078758616e-15"
$ws.Range("O5").Value = [double]"-0.1382206747322256"
$ws.Range("Q5").Value = [double]"-8.09642309068962e-16"
$ws.Range("R5").Value = [double]"-0.1237577023748392"
$ws.Range("S5").Value = [double]"0.03078035019887206"
$ws.Range("T5").Value = [double]"0.2052795968132872"
$ws.Range("U5").Value = [double]"-0.1686253156223198"
$ws.Range("B6").Value = [double]"0.0684526912477503"
$ws.Range("C6").Value = [double]"0.02012647564497602"
$ws.Range("D6").Value = [double]"0.5543526031480206"
$ws.Range("E6").Value = [double]"0.9818426010967805"
$ws.Range("G6").Value = [double]"0.1495776880695702"
$ws.Range("H6").Value = [double]"0.695923559201002"
$ws.Range("I6").Value = [double]"0.6536066037859465"
$ws.Range("J6").Value = [double]"0.7760011132205292"
$ws.Range("K6").Value = [double]"0.6158020315361707"
$ws.Range("L6").Value = [double]"0.2450391849937274"
$ws.Range("M6").Value = [double]"0.4535785045770876"
$ws.Range("N6").Value = [double]"1.862429420962686e-16"
$ws.Range("O6").Value = [double]"-0.1468461572647997"
$ws.Range("Q6").Value = [double]"2.058474623169284e-16"
$ws.Range("R6").Value = [double]"-0.1241251167351905"
$ws.Range("S6").Value = [double]"0.03942617309147878"
$ws.Range("T6").Value = [double]"0.2077386794830388"
$ws.Range("U6").Value = [double]"-0.1534767529571634"
$ws.Range("B7").Value = [double]"-0.03158809536073708"
$ws.Range("C7").Value = [double]"0.1693011165915949"
$ws.Range("D7").Value = [double]"0.5482195976645435"
$ws.Range("E7").Value = [double]"0.1320648805310834"
$ws.Range("F7").Value = [double]"0.1495776880695702"
$ws.Range("H7").Value = [double]"-0.0462145601951886"
$ws.Range("I7").Value = [double]"0.4973878604985454"
$ws.Range("J7").Value = [double]"0.05989144877916316"
$ws.Range("K7").Value = [double]"0.06057495152870172"
$ws.Range("L7").Value = [double]"0.925931183991245"
$ws.Range("M7").Value = [double]"0.01447525533394661"
$ws.Range("N7").Value = [double]"5.989274933711637e-16"
$ws.Range("O7").Value = [double]"0.07084167482301552"
$ws.Range("Q7").Value = [double]"6.316813406648991e-17"
$ws.Range("R7").Value = [double]"-0.07088220208111751"
$ws.Range("S7").Value = [double]"0.5452259023589349"
$ws.Range("T7").Value = [double]"0.1261931015454739"
$ws.Range("U7").Value = [double]"0.3200150253660413"
$ws.Range("B8").Value = [double]"0.003923958723663769"
$ws.Range("C8").Value = [double]"0.08724973029280536"
$ws.Range("D8").Value = [double]"0.3021693211650757"
$ws.Range("E8").Value = [double]"0.6773768707777231"
$ws.Range("F8").Value = [double]"0.695923559201002"
$ws.Range("G8").Value = [double]"-0.0462145601951886"
$ws.Range("I8").Value = [double]"0.4417259637869415"
$ws.Range("J8").Value = [double]"0.6783259659647563"
$ws.Range("K8").Value = [double]"0.5371519601751835"
$ws.Range("L8").Value = [double]"0.01902999694893956"
$ws.Range("M8").Value = [double]"0.6644594252512607"
$ws.Range("N8").Value = [double]"-6.393530636555416e-16"
$ws.Range("O8").Value = [double]"-0.1561372401110707"
$ws.Range("Q8").Value = [double]"-8.47776529340608e-16"
$ws.Range("R8").Value = [double]"0.07406909651085969"
$ws.Range("S8").Value = [double]"-0.1945245536176145"
$ws.Range("T8").Value = [double]"0.0955755836976475"
$ws.Range("U8").Value = [double]"-0.2878741473018495"
$ws.Range("B9").Value = [double]"0.03903941388948123"
$ws.Range("C9").Value = [double]"-0.05039689115876729"
$ws.Range("D9").Value = [double]"0.8662642689568739"
$ws.Range("E9").Value = [double]"0.6355324873451915"
$ws.Range("F9").Value = [double]"0.6536066037859465"
$ws.Range("G9").Value = [double]"0.4973878604985454"
$ws.Range("H9").Value = [double]"0.4417259637869415"
$ws.Range("J9").Value = [double]"0.6480552706311556"
$ws.Range("K9").Value = [double]"0.5231913527924974"
$ws.Range("L9").Value = [double]"0.5388699913355584"
$ws.Range("M9").Value = [double]"0.4011784428657405"
$ws.Range("N9").Value = [double]"-3.701870849802112e-17"
$ws.Range("O9").Value = [double]"-0.01339290751790123"
$ws.Range("Q9").Value = [double]"-5.182619189722956e-17"
$ws.Range("R9").Value = [double]"-0.2114340313067629"
$ws.Range("S9").Value = [double]"0.3863165654041744"
$ws.Range("T9").Value = [double]"0.2625638264202566"
$ws.Range("U9").Value = [double]"0.1777731337632806"
$ws.Range("B10").Value = [double]"0.1721239418668228"
$ws.Range("C10").Value = [double]"-0.08104252506935648"
$ws.Range("D10").Value = [double]"0.5044893085154387"
$ws.Range("E10").Value = [double]"0.7614718058408736"
$ws.Range("F10").Value = [double]"0.7760011132205292"
$ws.Range("G10").Value = [double]"0.05989144877916316"
$ws.Range("H10").Value = [double]"0.6783259659647563"
$ws.Range("I10").Value = [double]"0.6480552706311556"
$ws.Range("K10").Value = [double]"0.9443823590156931"
$ws.Range("L10").Value = [double]"0.129153649630667"
$ws.Range("M10").Value = [double]"0.6756281915557604"
$ws.Range("N10").Value = [double]"-1.701668213155127e-15"
$ws.Range("O10").Value = [double]"-0.1758183260735698"
$ws.Range("Q10").Value = [double]"-1.727647880531541e-15"
$ws.Range("R10").Value = [double]"-0.1858152972750236"
$ws.Range("S10").Value = [double]"-0.03366534267093271"
$ws.Range("T10").Value = [double]"0.1420392012739939"
$ws.Range("U10").Value = [double]"-0.1375284691121885"
$ws.Range("B11").Value = [double]"0.198667186370348"
$ws.Range("C11").Value = [double]"-0.05592584816530667"
$ws.Range("D11").Value = [double]"0.4317035091799344"
$ws.Range("E11").Value = [double]"0.6007500003764545"
$ws.Range("F11").Value = [double]"0.6158020315361707"
$ws.Range("G11").Value = [double]"0.06057495152870172"
$ws.Range("H11").Value = [double]"0.5371519601751835"
$ws.Range("I11").Value = [double]"0.5231913527924974"
$ws.Range("J11").Value = [double]"0.9443823590156931"
$ws.Range("L11").Value = [double]"0.1016463147940526"
$ws.Range("M11").Value = [double]"0.6966916537443587"
$ws.Range("N11").Value = [double]"-1.081275398247217e-15"
$ws.Range("O11").Value = [double]"-0.167532634084244"
$ws.Range("Q11").Value = [double]"-1.081275398247217e-15"
$ws.Range("R11").Value = [double]"-0.1920470706759466"
$ws.Range("S11").Value = [double]"-0.02977837223206488"
$ws.Range("T11").Value = [double]"0.06501534058395837"
$ws.Range("U11").Value = [double]"-0.09117699200325628"
$ws.Range("B12").Value = [double]"-0.05612525113176752"
$ws.Range("C12").Value = [double]"0.09867565890553331"
$ws.Range("D12").Value = [double]"0.5059907001350852"
$ws.Range("E12").Value = [double]"0.2210786576314376"
$ws.Range("F12").Value = [double]"0.2450391849937274"
$ws.Range("G12").Value = [double]"0.925931183991245"
$ws.Range("H12").Value = [double]"0.01902999694893956"
$ws.Range("I12").Value = [double]"0.5388699913355584"
$ws.Range("J12").Value = [double]"0.129153649630667"
$ws.Range("K12").Value = [double]"0.1016463147940526"
$ws.Range("M12").Value = [double]"-0.02959875824881635"
$ws.Range("N12").Value = [double]"2.278159390997511e-16"
$ws.Range("O12").Value = [double]"0.08304810588742481"
$ws.Range("Q12").Value = [double]"3.027143300366556e-16"
$ws.Range("R12").Value = [double]"-0.05417198617910895"
$ws.Range("S12").Value = [double]"0.4639179598906222"
$ws.Range("T12").Value = [double]"0.1508271710080992"
$ws.Range("U12").Value = [double]"0.2192401947192325"
$ws.Range("B13").Value = [double]"0.1130601880950822"
$ws.Range("C13").Value = [double]"0.1012401029327586"
$ws.Range("D13").Value = [double]"0.3384409619207953"
$ws.Range("E13").Value = [double]"0.4575882349190203"
$ws.Range("F13").Value = [double]"0.4535785045770876"
$ws.Range("G13").Value = [double]"0.01447525533394661"
$ws.Range("H13").Value = [double]"0.6644594252512607"
$ws.Range("I13").Value = [double]"0.4011784428657405"
$ws.Range("J13").Value = [double]"0.6756281915557604"
$ws.Range("K13").Value = [double]"0.6966916537443587"
$ws.Range("L13").Value = [double]"-0.02959875824881635"
$ws.Range("N13").Value = [double]"6.223814001381406e-16"
$ws.Range("O13").Value = [double]"-0.04109472597421545"
$ws.Range("Q13").Value = [double]"6.223814001381406e-16"
$ws.Range("R13").Value = [double]"-0.03380373102892087"
$ws.Range("S13").Value = [double]"0.02337922207127505"
$ws.Range("T13").Value = [double]"0.04573169243727224"
$ws.Range("U13").Value = [double]"-0.06029795865453617"
$ws.Range("D14").Value = [double]"7.426571720754092e-16"
$ws.Range("E14").Value = [double]"-1.079523078758616e-15"
$ws.Range("F14").Value = [double]"1.862429420962686e-16"
$ws.Range("G14").Value = [double]"5.989274933711637e-16"
$ws.Range("H14").Value = [double]"-6.393530636555416e-16"
$ws.Range("I14").Value = [double]"-3.701870849802112e-17"
$ws.Range("J14").Value = [double]"-1.701668213155127e-15"
$ws.Range("K14").Value = [double]"-1.081275398247217e-15"
$ws.Range("L14").Value = [double]"2.278159390997511e-16"
$ws.Range("M14").Value = [double]"6.223814001381406e-16"
$ws.Range("D15").Value = [double]"0.02970424846317648"
$ws.Range("E15").Value = [double]"-0.1382206747322256"
$ws.Range("F15").Value = [double]"-0.1468461572647997"
$ws.Range("G15").Value = [double]"0.07084167482301552"
$ws.Range("H15").Value = [double]"-0.1561372401110707"
$ws.Range("I15").Value = [double]"-0.01339290751790123"
$ws.Range("J15").Value = [double]"-0.1758183260735698"
$ws.Range("K15").Value = [double]"-0.167532634084244"
$ws.Range("L15").Value = [double]"0.08304810588742481"
$ws.Range("M15").Value = [double]"-0.04109472597421545"
$ws.Range("D17").Value = [double]"4.835907167002664e-16"
$ws.Range("E17").Value = [double]"-8.09642309068962e-16"
$ws.Range("F17").Value = [double]"2.058474623169284e-16"
$ws.Range("G17").Value = [double]"6.316813406648991e-17"
$ws.Range("H17").Value = [double]"-8.47776529340608e-16"
$ws.Range("I17").Value = [double]"-5.182619189722956e-17"
$ws.Range("J17").Value = [double]"-1.727647880531541e-15"
$ws.Range("K17").Value = [double]"-1.081275398247217e-15"
$ws.Range("L17").Value = [double]"3.027143300366556e-16"
$ws.Range("M17").Value = [double]"6.223814001381406e-16"
$ws.Range("D18").Value = [double]"-0.1987608672475611"
$ws.Range("E18").Value = [double]"-0.1237577023748392"
$ws.Range("F18").Value = [double]"-0.1241251167351905"
$ws.Range("G18").Value = [double]"-0.07088220208111751"
$ws.Range("H18").Value = [double]"0.07406909651085969"
$ws.Range("I18").Value = [double]"-0.2114340313067629"
$ws.Range("J18").Value = [double]"-0.1858152972750236"
$ws.Range("K18").Value = [double]"-0.1920470706759466"
$ws.Range("L18").Value = [double]"-0.05417198617910895"
$ws.Range("M18").Value = [double]"-0.03380373102892087"
$ws.Range("D19").Value = [double]"0.4302309024326345"
$ws.Range("E19").Value = [double]"0.03078035019887206"
$ws.Range("F19").Value = [double]"0.03942617309147878"
$ws.Range("G19").Value = [double]"0.5452259023589349"
$ws.Range("H19").Value = [double]"-0.1945245536176145"
$ws.Range("I19").Value = [double]"0.3863165654041744"
$ws.Range("J19").Value = [double]"-0.03366534267093271"
$ws.Range("K19").Value = [double]"-0.02977837223206488"
$ws.Range("L19").Value = [double]"0.4639179598906222"
$ws.Range("M19").Value = [double]"0.02337922207127505"
$ws.Range("D20").Value = [double]"0.2205306032889259"
$ws.Range("E20").Value = [double]"0.2052795968132872"
$ws.Range("F20").Value = [double]"0.2077386794830388"
$ws.Range("G20").Value = [double]"0.1261931015454739"
$ws.Range("H20").Value = [double]"0.0955755836976475"
$ws.Range("I20").Value = [double]"0.2625638264202566"
$ws.Range("J20").Value = [double]"0.1420392012739939"
$ws.Range("K20").Value = [double]"0.06501534058395837"
$ws.Range("L20").Value = [double]"0.1508271710080992"
$ws.Range("M20").Value = [double]"0.04573169243727224"
$ws.Range("D21").Value = [double]"0.2244033592937769"
$ws.Range("E21").Value = [double]"-0.1686253156223198"
$ws.Range("F21").Value = [double]"-0.1534767529571634"
$ws.Range("G21").Value = [double]"0.3200150253660413"
$ws.Range("H21").Value = [double]"-0.2878741473018495"
$ws.Range("I21").Value = [double]"0.1777731337632806"
$ws.Range("J21").Value = [double]"-0.1375284691121885"
$ws.Range("K21").Value = [double]"-0.09117699200325628"
$ws.Range("L21").Value = [double]"0.2192401947192325"
$ws.Range("M21").Value = [double]"-0.06029795865453617"
